$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '29.120.26'
    'E2' = '  +2.86%  '
    'D3' = '1.579.83'
    'E3' = '  +1.81%  '
    'D4' = '0.998'
    'E4' = '  -0.18%  '
    'D5' = '212.24'
    'E5' = '  +1.14%  '
    'D6' = '0.513'
    'E6' = '  +6.17%  '
    'E7' = '  -0.18%  '
    'D8' = '26.36'
    'E8' = '  +11.08%  '
    'E9' = '  +2.43%  '
    'D10' = '0.0593'
    'E10' = '  +1.64%  '
    'D11' = '0.0904'
    'D12' = '1.803.54'
    'E12' = '  +1.64%  '
    'D13' = '1.565.38'
    'E13' = '  +0.39%  '
    'D14' = '29.135.88'
    'E14' = '  +3.01%  '
    'D15' = '0.524'
    'E15' = '  +2.90%  '
    'E16' = '  +2.70%  '
    'D17' = '62.30'
    'E17' = '  +2.90%  '
    'D18' = '236.19'
    'E18' = '  +3.66%  '
    'E19' = '  +1.60%  '
    'E20' = '  +2.25%  '
    'D21' = '0.999'
    'E21' = '  -0.08%  '
    'D22' = '3.99'
    'E22' = '  +1.85%  '
    'D23' = '9.18'
    'E23' = '  +3.24%  '
    'E24' = '  +4.45%  '
    'D25' = '153.54'
    'E25' = '  +1.39%  '
    'D26' = '0.108'
    'E26' = '  +4.25%  '
    'D27' = '15.13'
    'E27' = '  +2.57%  '
    'D28' = '6.35'
    'E28' = '  +1.63%  '
    'D29' = '0.998'
    'E29' = '  -0.15%  '
    'D30' = '0.0468'
    'E30' = '  +0.24%  '
    'E31' = '  +0.27%  '
    'E32' = '  +1.50%  '
    'D33' = '1.422.47'
    'E33' = '  +2.48%  '
    'D34' = '3.07'
    'E34' = '  +1.24%  '
    'E35' = '  -1.58%  '
    'E36' = '  +1.74%  '
    'D37' = '2.74'
    'E37' = '  +5.76%  '
    'D38' = '2.30'
    'E38' = '  -1.71%  '
    'E39' = '  +1.27%  '
    'D40' = '0.530'
    'E40' = '  +3.83%  '
    'D41' = '1.97'
    'E41' = '  +1.87%  '
    'E42' = '  -0.10%  '
    'D43' = '52.71'
    'E43' = '  +24.49%  '
    'D44' = '0.788'
    'E44' = '  +1.28%  '
    'E45' = '  +0.38%  '
    'D46' = '64.64'
    'E46' = '  +4.40%  '
    'E47' = '  -0.82%  '
    'D48' = '1.715.78'
    'E48' = '  +1.60%  '
    'D49' = '0.848'
    'E49' = '  -6.37%  '
    'D50' = '85.66'
    'E50' = '  -0.09%  '
    'E51' = '  +1.61%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.ClearFormats()
}
